$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Diseases")

# Clear the "eyes" rule row (row 13) contents, keep formatting/style
$ws.Range("A13:C13").ClearContents()

# Update selection to reflect the new active cell/selection
$ws.Range("A13:C13").Select()

# Update window position/size
$excel.Left = 4340
$excel.Top = 500
$excel.Width = 24460
$excel.Height = 16400
